$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters that receive new values, in order, for each data row (rows 2-25)
$cols = @("B", "C", "D", "E", "G", "I", "J", "L", "O")

# New values per row (row 2 through row 25), one array per row matching $cols order
$data = @(
    @(1.921926567510866, 0.1576293812804295, 0.4733464799395364, 0.1420343467962297, 0.00251497375289989, 1.20108323770453, 0.04967698208863602, 0.5293025237902924, 5.505646947653986),
    @(1.794413545715486, 0.1381206173816167, 0.4729019006261979, 0.1429773773369547, 0.002518514085820055, 1.216236772433053, 0.04968975813735277, 0.519625658981127, 5.545406810180566),
    @(1.716419265835668, 0.1260893837454944, 0.4728408817741183, 0.1436088988382647, 0.002520804042656004, 1.226277287875238, 0.04970284770466904, 0.5138865691823895, 5.573292302903809),
    @(1.684713104324487, 0.1211735121667061, 0.4728693703399074, 0.1438794642519863, 0.002521766523011456, 1.230553806110898, 0.04970950752608516, 0.5115989371522716, 5.585528180451206),
    @(1.679453035790289, 0.1203564541497428, 0.4728773249854896, 0.143925190029881, 0.002521928114671581, 1.231275082179849, 0.04971069365589997, 0.5112221683467197, 5.587612596527066),
    @(1.715991350584488, 0.1260231391053424, 0.4728410498749014, 0.1436124942506058, 0.002520816904210333, 1.226334213994836, 0.04970293214462451, 0.5138555102830225, 5.573453789538888),
    @(1.877899062007828, 0.1509138404652219, 0.4731492303181, 0.1423486166663945, 0.002516170401284871, 1.206155266491784, 0.0496803027074435, 0.5259239954725956, 5.518634800008641),
    @(2.197707922326458, 0.1992982635304088, 0.4754334405441938, 0.14028607088607, 0.002507976319254071, 1.172434189336574, 0.04967725681518154, 0.551191988993196, 5.43873542406601),
    @(2.434015746378122, 0.2345800769569735, 0.4781340665134621, 0.139023447641927, 0.002502509795718057, 1.151236959580885, 0.04969984862867349, 0.5707282295802827, 5.396926733022525),
    @(2.541798784715638, 0.2505717871381989, 0.4795843724172641, 0.1385037547318824, 0.002500141924601415, 1.142373124848952, 0.04971544230854086, 0.5798260161096636, 5.381590001715438),
    @(2.582652988333848, 0.2566189023571894, 0.4801654178422297, 0.1383148101647116, 0.002499262274668967, 1.139128846052543, 0.04972210532835497, 0.5833012858962974, 5.376313052343903),
    @(2.573852595934511, 0.2553169340091301, 0.4800388636671045, 0.138355153721136, 0.002499450967746791, 1.139822562471025, 0.04972063671643312, 0.5825514865496046, 5.377425913405887),
    @(2.545159111897817, 0.2510694606020536, 0.4796315373898068, 0.1384880528492491, 0.002500069214786297, 1.142103964673176, 0.04971597532376393, 0.5801113255742791, 5.381145220056567),
    @(2.527588560123831, 0.2484666338343402, 0.4793861842442197, 0.138570479616801, 0.002500450122012167, 1.143516017068187, 0.04971321861352074, 0.5786205764468093, 5.383492557197826),
    @(2.426977448615844, 0.2335337905606423, 0.4780437444600807, 0.1390585101553068, 0.002502666926315397, 1.151831927472927, 0.04969893597222708, 0.570137893315021, 5.398003226283947),
    @(2.365327508016435, 0.2243579098376074, 0.4772769713039935, 0.1393718981572771, 0.002504057252437422, 1.157133146307899, 0.04969153109516888, 0.5649878834715309, 5.407849042529733),
    @(2.32989507244173, 0.21907472510037, 0.4768568211003554, 0.139557298118314, 0.002504868125504616, 1.160255564551662, 0.04968777280757841, 0.5620455708648819, 5.413858612587461),
    @(2.317902951798658, 0.2172849988226346, 0.4767181525686652, 0.1396209557844834, 0.002505144598481945, 1.161325342727309, 0.04968658652763835, 0.5610527661577436, 5.41595282682357),
    @(2.371887471255377, 0.225335264892351, 0.4773564352093729, 0.1393380048066568, 0.002503908091986741, 1.15656123462697, 0.04969226757332024, 0.5655340589104298, 5.406765066420377),
    @(2.553586037686387, 0.2523172817118962, 0.4797503150563927, 0.1384488041451277, 0.002499887159604778, 1.141430812686679, 0.04971732396664841, 0.5808272433697823, 5.380038355608662),
    @(2.672563613786053, 0.2699013652253655, 0.4815004564524088, 0.1379134216431073, 0.002497358365236946, 1.132196683107544, 0.04973811520194715, 0.5909978050272429, 5.365665026682763),
    @(2.609042900044017, 0.2605210810410483, 0.4805494052356352, 0.1381949818399661, 0.002498698988368832, 1.137065141389527, 0.04972661663819444, 0.5855535702153674, 5.37305282194373),
    @(2.368921677662343, 0.2248934271066219, 0.4773204451668107, 0.1393533117016901, 0.002503975491451137, 1.156819563295613, 0.04969193305699804, 0.5652870753697385, 5.407254044720872),
    @(2.110950635140227, 0.1862552960271842, 0.4746358349267439, 0.1407996031601506, 0.002510095389949525, 1.18092919267805, 0.04967369427326496, 0.5441853352584474, 5.45738854373127)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $colLetter = $cols[$j]
        $cellRef = "$colLetter$rowNum"
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}
